$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, and the percentage-volume column).
# None of these values look like numbers to Excel, so a normal .Value
# assignment keeps them as text without any extra handling.
$textUpdates = @{
    "E2" = "  +3.87%  "
    "E3" = "  +0.84%  "
    "E4" = "  +0.24%  "
    "E5" = "  -0.14%  "
    "E6" = "  +0.14%  "
    "E7" = "  +1.75%  "
    "E8" = "  +2.52%  "
    "B9" = "Dogecoin"
    "C9" = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
    "E9" = "  +1.52%  "
    "B10" = "OKB"
    "C10" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "E10" = "  -1.08%  "
    "E11" = "  +1.73%  "
    "E12" = "  +0.15%  "
    "E13" = "  +0.86%  "
    "E14" = "  +1.13%  "
    "E15" = "  +5.52%  "
    "E16" = "  +0.97%  "
    "E17" = "  +0.61%  "
    "E18" = "  +0.39%  "
    "E19" = "  +0.03%  "
    "E20" = "  +0.13%  "
    "E21" = "  +2.61%  "
    "E22" = "  +0.71%  "
    "E23" = "  +3.58%  "
    "E24" = "  +0.94%  "
    "E25" = "  +1.38%  "
    "E26" = "  +3.49%  "
    "E27" = "  +1.75%  "
    "E28" = "  +1.81%  "
    "E29" = "  +0.21%  "
    "E30" = "  +2.02%  "
    "E31" = "  +2.84%  "
    "E32" = "  +1.28%  "
    "E33" = "  +1.88%  "
    "E34" = "  +1.27%  "
    "E35" = "  +12.36%  "
    "E36" = "  +9.56%  "
    "E37" = "  +2.78%  "
    "E38" = "  +4.01%  "
    "E39" = "  +5.80%  "
    "E40" = "  +1.72%  "
    "E41" = "  +3.14%  "
    "E42" = "  +7.07%  "
    "E43" = "  +0.10%  "
    "E44" = "  -2.25%  "
    "E45" = "  +3.22%  "
    "E46" = "  +2.52%  "
    "E47" = "  +1.37%  "
    "E48" = "  +2.60%  "
    "E49" = "  +4.82%  "
    "E50" = "  +4.57%  "
    "E51" = "  +3.17%  "
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Price column updates. These strings look numeric to Excel (e.g. "1.003",
# "0.07510", "28.404.61"), so a plain .Value assignment risks Excel silently
# reinterpreting them as real numbers -- which would drop meaningful trailing
# zeros, flip tiny values into scientific notation, and change the cell away
# from a text type entirely. To keep them as the exact literal text from the
# source data, temporarily force a text number format before assigning the
# value, then restore the cell's original style so no visible formatting
# change (or stray style) is left behind.
$priceUpdates = @{
    "D2" = "28.404.61"
    "D3" = "1.794.71"
    "D4" = "1.003"
    "D5" = "313.58"
    "D6" = "1.002"
    "D7" = "0.5311"
    "D8" = "0.3788"
    "D9" = "0.07510"
    "D10" = "42.50"
    "D11" = "1.111"
    "D13" = "20.95"
    "D14" = "6.171"
    "D15" = "7.362"
    "D16" = "1.793.96"
    "D17" = "90.17"
    "D18" = "0.00001064"
    "D19" = "0.06439"
    "D20" = "1.001"
    "D21" = "17.20"
    "D22" = "5.916"
    "D23" = "28.365.66"
    "D24" = "11.33"
    "D25" = "2.131"
    "D26" = "160.64"
    "D27" = "20.53"
    "D28" = "2.383"
    "D29" = "1.986.95"
    "D30" = "123.22"
    "D31" = "1.122"
    "D32" = "0.1011"
    "D33" = "5.696"
    "D34" = "3.669"
    "D35" = "0.2300"
    "D36" = "0.06536"
    "D37" = "0.02313"
    "D38" = "5.088"
    "D39" = "8.653"
    "D40" = "11.52"
    "D41" = "0.6323"
    "D42" = "1.210"
    "D43" = "1.001"
    "D44" = "1.403"
    "D45" = "13.55"
    "D46" = "0.5930"
    "D47" = "3.672"
    "D48" = "124.76"
    "D49" = "1.986"
    "D50" = "1.163"
    "D51" = "0.06937"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = $origStyle
}
